$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47 (shifts rows 47..136 down to 48..137)
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record
$ws.Cells.Item(47, 1).Value = 11
$ws.Cells.Item(47, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(47, 3).Value = "Bíobío"
$ws.Cells.Item(47, 4).Value = 44965
$ws.Cells.Item(47, 5).Value = 8
$ws.Cells.Item(47, 6).Value = "Fruta"
$ws.Cells.Item(47, 7).Value = 100101
$ws.Cells.Item(47, 8).Value = "Berries"
$ws.Cells.Item(47, 9).Value = 100101001
$ws.Cells.Item(47, 10).Value = "Arándano (blue)"
$ws.Cells.Item(47, 11).Value = "Sin especificar"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 200
$ws.Cells.Item(47, 14).Value = 3000
$ws.Cells.Item(47, 15).Value = 3500
$ws.Cells.Item(47, 16).Value = 3250
$ws.Cells.Item(47, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(47, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(47, 19).Value = 1625
$ws.Cells.Item(47, 20).Value = 2
